$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> new DAMSLTag (col I / 9), new DialogAct (col J / 10)
$updates = @(
    ,@(15, 'b', 'Acknowledge (Backchannel)')
    ,@(17, 'sd', 'Statement-non-opinion')
    ,@(18, 'b', 'Acknowledge (Backchannel)')
    ,@(20, 'sv', 'Statement-opinion')
    ,@(21, 'aa', 'Agree/Accept')
    ,@(27, 'b', 'Acknowledge (Backchannel)')
    ,@(28, 'sd', 'Statement-non-opinion')
    ,@(30, 'b', 'Acknowledge (Backchannel)')
    ,@(39, 'sd', 'Statement-non-opinion')
    ,@(40, 'sd', 'Statement-non-opinion')
    ,@(57, 'sd', 'Statement-non-opinion')
    ,@(69, 'sv', 'Statement-opinion')
    ,@(72, 'sd', 'Statement-non-opinion')
    ,@(73, 'aa', 'Agree/Accept')
    ,@(90, 'sd', 'Statement-non-opinion')
    ,@(94, 'sd', 'Statement-non-opinion')
    ,@(108, 'sd', 'Statement-non-opinion')
    ,@(129, 'sv', 'Statement-opinion')
    ,@(134, 'b', 'Acknowledge (Backchannel)')
    ,@(138, 'sd', 'Statement-non-opinion')
    ,@(141, '%', 'Uninterpretable')
    ,@(142, '%', 'Uninterpretable')
    ,@(146, 'sd', 'Statement-non-opinion')
    ,@(154, '%', 'Uninterpretable')
    ,@(166, 'ba', 'Appreciation')
    ,@(172, 'aa', 'Agree/Accept')
    ,@(173, 'b', 'Acknowledge (Backchannel)')
    ,@(175, 'ba', 'Appreciation')
    ,@(184, 'sd', 'Statement-non-opinion')
    ,@(185, 'aa', 'Agree/Accept')
    ,@(187, 'sd', 'Statement-non-opinion')
    ,@(196, 'b', 'Acknowledge (Backchannel)')
    ,@(215, 'sv', 'Statement-opinion')
    ,@(216, 'sd', 'Statement-non-opinion')
    ,@(241, 'aa', 'Agree/Accept')
    ,@(244, 'ba', 'Appreciation')
    ,@(248, 'b', 'Acknowledge (Backchannel)')
    ,@(256, 'b', 'Acknowledge (Backchannel)')
    ,@(258, 'sv', 'Statement-opinion')
    ,@(261, 'aa', 'Agree/Accept')
    ,@(275, 'sd', 'Statement-non-opinion')
    ,@(286, '%', 'Uninterpretable')
    ,@(288, 'sd', 'Statement-non-opinion')
    ,@(289, 'sd', 'Statement-non-opinion')
    ,@(301, 'sd', 'Statement-non-opinion')
    ,@(302, 'sd', 'Statement-non-opinion')
    ,@(331, '%', 'Uninterpretable')
    ,@(336, 'sd', 'Statement-non-opinion')
    ,@(340, 'b', 'Acknowledge (Backchannel)')
    ,@(341, 'sv', 'Statement-opinion')
    ,@(345, 'sd', 'Statement-non-opinion')
    ,@(346, 'sd', 'Statement-non-opinion')
    ,@(349, 'aa', 'Agree/Accept')
    ,@(352, 'sd', 'Statement-non-opinion')
    ,@(356, 'sd', 'Statement-non-opinion')
    ,@(359, '%', 'Uninterpretable')
    ,@(361, 'sv', 'Statement-opinion')
    ,@(374, 'b', 'Acknowledge (Backchannel)')
    ,@(382, 'sv', 'Statement-opinion')
    ,@(397, 'sv', 'Statement-opinion')
    ,@(405, 'aa', 'Agree/Accept')
    ,@(408, 'aa', 'Agree/Accept')
    ,@(410, '%', 'Uninterpretable')
    ,@(416, 'ba', 'Appreciation')
    ,@(425, 'ba', 'Appreciation')
    ,@(432, 'aa', 'Agree/Accept')
    ,@(442, 'sv', 'Statement-opinion')
    ,@(452, 'b', 'Acknowledge (Backchannel)')
    ,@(457, 'ba', 'Appreciation')
    ,@(460, 'qy', 'Yes-No-Question')
    ,@(465, 'sd', 'Statement-non-opinion')
    ,@(468, 'aa', 'Agree/Accept')
    ,@(470, 'b', 'Acknowledge (Backchannel)')
    ,@(487, 'ba', 'Appreciation')
    ,@(493, 'aa', 'Agree/Accept')
    ,@(494, 'b', 'Acknowledge (Backchannel)')
    ,@(496, 'b', 'Acknowledge (Backchannel)')
    ,@(499, 'sv', 'Statement-opinion')
    ,@(500, 'sd', 'Statement-non-opinion')
    ,@(501, 'sv', 'Statement-opinion')
    ,@(506, 'sv', 'Statement-opinion')
    ,@(515, 'sd', 'Statement-non-opinion')
    ,@(524, 'b', 'Acknowledge (Backchannel)')
    ,@(526, 'sd', 'Statement-non-opinion')
    ,@(528, 'sd', 'Statement-non-opinion')
    ,@(529, 'sd', 'Statement-non-opinion')
    ,@(530, 'b', 'Acknowledge (Backchannel)')
    ,@(536, 'sv', 'Statement-opinion')
    ,@(546, 'b', 'Acknowledge (Backchannel)')
    ,@(574, 'sv', 'Statement-opinion')
    ,@(585, 'aa', 'Agree/Accept')
    ,@(604, 'ba', 'Appreciation')
    ,@(651, 'b', 'Acknowledge (Backchannel)')
    ,@(661, 'b', 'Acknowledge (Backchannel)')
    ,@(663, 'b', 'Acknowledge (Backchannel)')
    ,@(685, 'sd', 'Statement-non-opinion')
    ,@(726, 'sd', 'Statement-non-opinion')
    ,@(732, 'aa', 'Agree/Accept')
    ,@(743, 'sd', 'Statement-non-opinion')
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 9).Value = $u[1]
    $ws.Cells.Item($r, 10).Value = $u[2]
}
